$d = $word.ActiveDocument

# The document contains several repeated "1.Semalam saya.. / 2. Hari ini
# saya, / 3. Saya ada kendala " placeholder blocks. The edit only targets
# the block right under the final "Ricky :" heading near the end of the
# document (Word paragraph indices 313-315).
$p1 = $d.Paragraphs.Item(313)
$p1.Range.Find.Execute("1.Semalam saya.. ", $false, $false, $false, $false, $false, $true, 1, $false, "1.Semalam saya menambah efek ke semua slide.", 2)

$p2 = $d.Paragraphs.Item(314)
$p2.Range.Find.Execute("2. Hari ini saya, ", $false, $false, $false, $false, $false, $true, 1, $false, "2. Hari ini saya melakukan sedikit perubahan pada slide pengantar dan menambah note.", 2)

# Third paragraph: "3. Saya ada kendala " becomes "3. Saya" + " tidak" +
# " ada kendala." split across three runs.
$p3 = $d.Paragraphs.Item(315)
$p3.Range.Text = "3. Saya"

$p3 = $d.Paragraphs.Item(315)
$insertPos = $p3.Range.End - 1
$ip = $d.Range($insertPos, $insertPos)
$ip.InsertAfter(" tidak ada kendala.")

$p3 = $d.Paragraphs.Item(315)
$splitStart = $p3.Range.Start + 7
$splitEnd = $splitStart + 6
$midRange = $d.Range($splitStart, $splitEnd)
$midRange.Bold = 1
$midRange.Bold = 0

$tailStart = $splitEnd
$tailEnd = $p3.Range.End - 1
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Bold = 1
$tailRange.Bold = 0
